$p = $ppt.ActivePresentation

# Add a new slide at the end (index 17) using the "Title Only" layout,
# matching the other appended slide's minimal placeholder set.
$s = $p.Slides.Add(17, 11)

# Title placeholder: "Presented to you by"
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Presented to you by"

# Free-floating text box with the presenter's name.
# AddTextbox takes Left/Top/Width/Height in points (1 pt = 12700 EMU);
# target EMU: off (8333752,5251304) ext (3633165,1325563).
$tb = $s.Shapes.AddTextbox(1, 656.20094488188975, 413.48850393700787, 286.07598425196850, 104.37503937007874)
$tb.TextFrame.TextRange.Text = "Dustin Walker"
